# Applies the "Updated symbol list" commit (Wed Dec 21 2022) to the
# cryptocurrency tracker sheet:
#   - Refreshed Price (column D) figures for most rows.
#   - Rows 41-43 (KickToken / BKEXToken / CEJI) were re-ranked, so their
#     Coin name (B), Link (C), Price (D) and rank label (E) values rotate
#     between the three rows.
#   - E19's rank label lost its "Worstin24h" suffix; E43 gained it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") holds numeric-looking text that must stay text ---
# Assigning a plain numeric-looking string to .Value makes Excel coerce it
# to a real number, which loses the exact source formatting/precision
# (e.g. trailing zeros such as "0.1070" or "0.003040", or rounds values
# like "22.72" into binary floating point noise). To keep these as text,
# each target cell is temporarily switched to the "@" (Text) number
# format before the value is written, then restored to the default
# "Normal" style afterwards so no visible formatting change remains.
$priceValues = [ordered]@{
    "D2" = "248.98"
    "D3" = "22.72"
    "D4" = "5.271"
    "D5" = "0.05704"
    "D6" = "3.405"
    "D7" = "6.321"
    "D8" = "0.8047"
    "D9" = "0.8945"
    "D10" = "0.1429"
    "D11" = "0.07417"
    "D14" = "0.09398"
    "D15" = "3.859"
    "D16" = "0.001578"
    "D17" = "0.04795"
    "D18" = "0.01827"
    "D19" = "0.0005796"
    "D20" = "0.006426"
    "D21" = "0.004994"
    "D22" = "0.0009964"
    "D23" = "0.0001499"
    "D24" = "3.698"
    "D25" = "2.198"
    "D27" = "0.1353"
    "D40" = "0.03978"
    "D41" = "0.1070"
    "D42" = "0.002729"
    "D43" = "0.003040"
    "D44" = "0.007691"
    "D45" = "0.00005567"
    "D46" = "0.00000000749"
    "D47" = "0.4986"
    "D48" = "0.2037"
    "D49" = "0.00002098"
    "D50" = "0.01009"
}
foreach ($addr in $priceValues.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $priceValues.Keys) {
    $ws.Range($addr).Value = $priceValues[$addr]
}
foreach ($addr in $priceValues.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# --- Columns B, C, E: plain text updates (no numeric coercion risk) ---
$ws.Range("E19").Value = "18OneONE"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
